$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C) for rows 2-8 from 45243 (2023-11-13) to 45244 (2023-11-14)
$ws.Range("C2:C8").Value = 45244
